# Fixed data errors for citations and language toggle #124
#
# 1) Two Spanish "Article"/"Press Release" rows on the
#    "pages_with_citations" sheet had their expectedHeaderText column
#    (D) filled with the English reference labels instead of the
#    correct Spanish ones. Correct them:
#      D4: "Selected References" -> "Bibliografía selecta"
#      D5: "Referencias"         -> "Bibliografía"
#
# 2) Language toggle: the active sheet/tab and the remembered
#    selection on each sheet were wrong (pages_without_citations was
#    shown/selected instead of pages_with_citations). Toggle the
#    active sheet back to "pages_with_citations" and restore the
#    expected selections on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("pages_with_citations")
$ws2 = $wb.Worksheets.Item("pages_without_citations")

# --- Fix the mislabeled citation references ---
$ws1.Range("D4").Value = "Bibliografía selecta"
$ws1.Range("D5").Value = "Bibliografía"

# --- Restore selection on the non-active sheet first ---
$ws2.Range("C12").Select()

# --- Toggle back to pages_with_citations as the active/selected tab ---
$ws1.Activate()
$ws1.Range("D4:D5").Select()
